# Fix typo'd values in column B (stray leading space duplicates of
# existing category labels) on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToFix = @(3, 8, 9, 11, 13, 15, 16, 19, 21, 23, 24)

foreach ($r in $rowsToFix) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2.Trim()
}

# Reflect the author's final selection/scroll position after editing.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B25").Select()
